# "Tuning commandes et winch"
# Updates two "(winch arm) position" textboxes: resize/reposition them and
# change their label text from placeholder strings to "Bras en position basse",
# split across 4 runs (to mirror the authored run-breaks "Bras "/"en"/" position "/"basse").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - TextBox 51 (id=90): "Brake swerve" -> "Bras en position basse"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(25)

# Resize (position stays the same). Shape.Width/Height round-trip through a
# 32-bit float, so the literals below are chosen to land exactly on the
# authored EMU values (1571625 x 275545) after that quantization.
$shp1.Width = 123.75
$shp1.Height = 21.6964566929134

$tr1 = $shp1.TextFrame.TextRange
$tr1.Text = "Bras en position basse"

# Re-apply the run-level font so the text splits into four runs, matching
# the authored "Bras " / "en" / " position " / "basse" run breaks.
$tr1.Characters(1, 5).Font.Name = "Calibri"
$tr1.Characters(6, 2).Font.Name = "Calibri"
$tr1.Characters(8, 10).Font.Name = "Calibri"
$tr1.Characters(18, 5).Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# Slide 2 - TextBox 47 (id=131): "(rien)" -> "Bras en position basse"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(31)

# Reposition + resize (see note above re: float32 quantization).
$shp2.Left = 315.954881889764
$shp2.Top = 78.859846519685
$shp2.Width = 123.75
$shp2.Height = 21.6964566929134

$tr2 = $shp2.TextFrame.TextRange

# Paragraph is no longer centered; it now uses the (default) left alignment.
$tr2.ParagraphFormat.Alignment = 1

$tr2.Text = "Bras en position basse"
$tr2.Characters(1, 5).Font.Name = "Calibri"
$tr2.Characters(6, 2).Font.Name = "Calibri"
$tr2.Characters(8, 10).Font.Name = "Calibri"
$tr2.Characters(18, 5).Font.Name = "Calibri"
